# Fixed naive component forecaster bug - Presentation state 11.02.
#
# For each data row (rows 2-16, corresponding to quarters Q0..Q9 of
# matched-error history), a newly computed error value is inserted into
# column B (the most-recent-quarter slot). The values that used to occupy
# columns B..J shift one column to the right (into C..K), and whatever
# value used to sit in the row's last used column falls off the end.
#
# New "column B" values introduced by this fix:
$newFirstValues = @{
    2  = 0.7496711949059137
    3  = -0.1976049264540507
    4  = 0.1255103924969555
    5  = -0.003050974015260888
    6  = 1.419652253737389
    7  = 0.1883110177716323
    8  = 0.2202779152847414
    9  = 0.4485660054549828
    10 = 0.420735823599318
    11 = -0.1252553916527783
    12 = 0.09764018641116785
    13 = -0.1898380159455487
    14 = 0.1743923273248104
    15 = -0.254916590923889
    16 = -0.01954473626955232
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B is column index 2. Data lives in columns B..K (2..11).
$firstCol = 2
$lastCol = 11

for ($row = 2; $row -le 16; $row++) {

    # Find the last populated column in this row within B..K, based on the
    # values already present before the edit.
    $lastUsedCol = $firstCol - 1
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($row, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val -ne "") {
            $lastUsedCol = $c
        }
    }

    # Shift existing values one column to the right, starting from the
    # rightmost populated column moving left, so we don't overwrite values
    # before they are read. The value in the last used column is dropped
    # if it would spill past column K.
    for ($c = $lastUsedCol; $c -ge $firstCol; $c--) {
        $srcVal = $ws.Cells.Item($row, $c).Value()
        $destCol = $c + 1
        if ($destCol -le $lastCol) {
            $ws.Cells.Item($row, $destCol).Value2 = $srcVal
        }
    }

    # Write the newly computed value into column B.
    $ws.Cells.Item($row, $firstCol).Value2 = $newFirstValues[$row]
}
